# Updates Jogos_da_Semana_FlashScore_2024-11-06 worksheet:
#  - tweak a handful of odds values in existing rows 5 and 7
#  - append four new match rows (8-11) with full odds data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Small numeric corrections on existing rows
# ---------------------------------------------------------------------------
$ws.Range("J5").Value  = 3.1
$ws.Range("N5").Value  = 7.5
$ws.Range("AR5").Value = 81

$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10

# ---------------------------------------------------------------------------
# 2. Append four new rows (8-11) of match data
# ---------------------------------------------------------------------------

# Column letters A..BD, in sheet order
$cols = @(
    "A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T",
    "U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL",
    "AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD"
)

# Row data, one array of 56 values (as text) per new row, in column order above.
$newRowsData = @(
    @("b9mfwU3n","06/11/2024","22:00","MEXICO - LIGA MX","Club America","Pachuca","1.7","3.9","4.75","2.3","2.3","5","1.04","13","1.25","3.75","1.75","2.05","1.33","3.25","1.8","1.95","7.5","8.5","8.5","13","13","23","12","7.5","15","51","13","23","15","51","34","41","201","3.75","8.5","19","29","41","126","3.25","8","51","6.5","23","29","81","101","201","501","126"),
    @("rcDpJQIU","06/11/2024","22:00","MEXICO - LIGA MX","Necaxa","Atl. San Luis","2.5","3.5","2.7","3","2.3","3.2","1.03","15","1.18","4.5","1.65","2.2","1.33","3.25","1.57","2.25","11","15","10","26","19","23","15","7","12","34","11","15","10","26","19","23","126","4.75","13","19","41","51","126","3.25","7","41","5","15","21","41","51","126","401","151"),
    @("lbKyLnlI","06/11/2024","22:00","MEXICO - LIGA MX","U.A.N.L.- Tigres","Toluca","1.73","3.9","4.5","2.25","2.5","4.33","1.03","17","1.17","5","1.53","2.4","1.25","3.75","1.53","2.38","10","10","8.5","15","12","19","17","8","12","34","17","26","15","41","29","29","126","4","8.5","15","26","41","81","3.75","7","41","6.5","21","23","67","67","126","351","151"),
    @("6ihZt8v3","06/11/2024","22:00","MEXICO - LIGA DE EXPANSION MX","Venados","Celaya","2.25","3.2","3","2.8","2.1","3.5","1.02","7.4","1.28","3.05","1.87","1.83","1.4","2.52","1.65","2","8.25","11.75","8.75","23","17.5","25","9.75","6.2","12.5","55","9.5","16","10.5","40","26","32","400","4.25","11.5","17.5","45","70","200","2.62","6.5","50","5","16.5","22","75","100","250","51","51")
)

$startRow = 8

for ($i = 0; $i -lt $newRowsData.Length; $i++) {
    $rowNum  = $startRow + $i
    $rowVals = $newRowsData[$i]

    for ($j = 0; $j -lt $cols.Length; $j++) {
        $colLetter = $cols[$j]
        $rawValue  = $rowVals[$j]
        $cell      = $ws.Range($colLetter + $rowNum)

        if ($j -lt 6) {
            # Text columns: Id, Date, Time, League, Home, Away.
            # The "Date" column (B) looks like a date (dd/mm/yyyy) and would
            # otherwise get auto-converted into a date serial number, so we
            # force it to stay as plain text first.
            if ($colLetter -eq "B") {
                $cell.NumberFormat = "@"
                $cell.Value = $rawValue
                $cell.Style = "Normal"
            } else {
                $cell.Value = $rawValue
            }
        } else {
            # Numeric odds columns
            $cell.Value = [double]$rawValue
        }
    }
}

Write-Output "Applied updates: tweaked rows 5 & 7, added rows 8-11"
